# Update the cryptocurrency prices/volumes for the daily refresh.
# Values that look numeric (e.g. "6.22", "0.999") must be written back as
# literal text (matching the sheet's existing inline-string cells), so each
# cell is forced to text format before the value is assigned and the
# number-format override is cleared again afterwards via the Normal style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# row -> @{ col letter = new value }
$updates = [ordered]@{
    2  = @{ D = '64.071.32';   E = '  +0.74%  ' }
    3  = @{ D = '2.757.08';    E = '  +1.12%  ' }
    4  = @{               E = '  -0.31%  ' }
    5  = @{ D = '577.75';      E = '  -1.16%  ' }
    6  = @{ D = '158.50';      E = '  +3.25%  ' }
    7  = @{               E = '  +0.13%  ' }
    8  = @{               E = '  +0.12%  ' }
    9  = @{               E = '  -1.48%  ' }
    10 = @{               E = '  -13.61%  ' }
    11 = @{ D = '0.388';       E = '  -0.78%  ' }
    12 = @{ D = '0.158';       E = '  -1.85%  ' }
    13 = @{ D = '3.245.94';    E = '  +0.49%  ' }
    14 = @{ D = '26.91';       E = '  +2.12%  ' }
    15 = @{ D = '63.957.99';   E = '  +0.68%  ' }
    16 = @{ D = '0.0000153';   E = '  +0.84%  ' }
    17 = @{ D = '2.761.15';    E = '  +0.35%  ' }
    18 = @{               E = '  +1.13%  ' }
    19 = @{ D = '4.90';        E = '  +0.34%  ' }
    20 = @{ D = '360.38';      E = '  -0.25%  ' }
    21 = @{               E = '  -2.12%  ' }
    22 = @{ D = '0.557';       E = '  +3.98%  ' }
    23 = @{ D = '0.999';       E = '  +0.19%  ' }
    24 = @{ D = '66.22';       E = '  +0.39%  ' }
    25 = @{ D = '0.171';       E = '  +2.03%  ' }
    26 = @{ D = '8.52';        E = '  -0.57%  ' }
    27 = @{               E = '  +0.06%  ' }
    28 = @{               E = '  +5.36%  ' }
    29 = @{ D = '1.97';        E = '  -2.62%  ' }
    30 = @{               E = '  -0.51%  ' }
    31 = @{               E = '  +3.32%  ' }
    32 = @{ D = '169.92';      E = '  -2.63%  ' }
    33 = @{ D = '20.42';       E = '  -0.55%  ' }
    34 = @{ D = '4.96';        E = '  +3.76%  ' }
    35 = @{               E = '  +0.00%  ' }
    36 = @{               E = '  +1.89%  ' }
    37 = @{ D = '1.81';        E = '  -0.43%  ' }
    38 = @{ D = '0.997';       E = '  -0.08%  ' }
    39 = @{ D = '6.22';        E = '  +12.06%  ' }
    40 = @{ D = '4.18';        E = '  -0.99%  ' }
    41 = @{ D = '331.25';      E = '  -3.94%  ' }
    42 = @{ D = '39.22';       E = '  +0.23%  ' }
    43 = @{ D = '21.86';       E = '  +0.61%  ' }
    44 = @{ D = '0.0596';      E = '  +1.49%  ' }
    45 = @{ D = '21.87';       E = '  -0.04%  ' }
    # Rows 46/47 swap their rank: Mantle moves up above VeChain.
    46 = @{ B = 'Mantle';  C = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt';          D = '0.637';  E = '  -1.42%  ' }
    47 = @{ B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet';      D = '0.0257'; E = '  +0.30%  ' }
    48 = @{ D = '136.44';      E = '  -2.34%  ' }
    49 = @{               E = '  +0.86%  ' }
    50 = @{               E = '  +0.37%  ' }
    51 = @{               E = '  +0.72%  ' }
}

$colNumbers = @{ A = 1; B = 2; C = 3; D = 4; E = 5 }

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($colLetter in $cols.Keys) {
        $value = $cols[$colLetter]
        $colNum = $colNumbers[$colLetter]
        if ($colLetter -eq 'D' -or $colLetter -eq 'E') {
            Set-TextValue $row $colNum $value
        } else {
            $ws.Cells.Item($row, $colNum).Value = $value
        }
    }
}
